$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.452
$ws.Range("D5").Value = 0.549
$ws.Range("E5").Value = 0.579
$ws.Range("F5").Value = 0.613
$ws.Range("G5").Value = 0.628
$ws.Range("H5").Value = 0.642

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.452
$ws.Range("E7").Value = 0.579
$ws.Range("F7").Value = 0.613

# Row 8: BERT-base
$ws.Range("C8").Value = 0.47
$ws.Range("D8").Value = 0.639
$ws.Range("E8").Value = 0.657
$ws.Range("F8").Value = 0.698
$ws.Range("G8").Value = 0.718
$ws.Range("H8").Value = 0.722

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.569
$ws.Range("C9").Value = 0.631
$ws.Range("D9").Value = 0.689
$ws.Range("E9").Value = 0.695
$ws.Range("F9").Value = 0.733
$ws.Range("G9").Value = 0.747
$ws.Range("H9").Value = 0.752
